# Automatische test-sync: 2025-07-22 12:18:50
# Adds a new test mail (row 5) to the "Logs" sheet and updates the
# "Dashboard" summary sheet (counts + ordering) to match.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs: append new row 5 -------------------------------------------------
$logs.Range("A5").Value = "Ik stuur het pakket terug."
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Testmail #4: Ik stuur het pakket terug."
$logs.Range("D5").Value = "Retour / Terugbetaling"
$logs.Range("E5").Value = "Beste klant,`nBedankt voor uw bericht. Om uw retourzending zo soepel mogelijk te laten verlopen, vragen wij u vriendelijk om het volgende te doen:`n- Vul het retourformulier in dat bij uw bestelling zat en voeg dit toe aan het pakket.`n- Stuur het pakket terug naar het volgende adres: [adres retourzending].`n- Zodra wij uw retourzending hebben ontvangen, zullen wij het verder afhandelen en u op de hoogte houden van de status van uw retour.`nMocht u nog verdere vragen of opmerkingen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Range("F5").Value = "2025-07-22 12:18:21"
$logs.Range("G5").Value = "Ja"
$logs.Range("H5").Value = "Nee"
$logs.Range("I5").Value = "Ja"
$logs.Range("J5").Value = "Ja"

# --- Logs: extend the existing conditional-formatting ranges to cover row 5
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "4")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "5")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard: reorder by updated counts (Retour now has 2 mails) --------
$dash.Range("A2").Value = "Retour / Terugbetaling"
$dash.Range("B2").Value = 2
$dash.Range("A3").Value = "Openingstijden / Locatie"
$dash.Range("B3").Value = 1
$dash.Range("A4").Value = "Bestelling / Levering"
$dash.Range("B4").Value = 1
